$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.899.50'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.413.06'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.50'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.07'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +6.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  +5.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.69'
$ws.Range('E11').Value = '  +1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000218'
$ws.Range('E12').Value = '  +47.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.14'
$ws.Range('E13').Value = '  +9.11%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.967.52'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.27'
$ws.Range('E16').Value = '  +7.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.403.09'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.45'
$ws.Range('E18').Value = '  +7.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.09'
$ws.Range('E19').Value = '  +6.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '61.945.63'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '458.69'
$ws.Range('E21').Value = '  +47.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '91.99'
$ws.Range('E22').Value = '  +8.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.21'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.07'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.28'
$ws.Range('E25').Value = '  +3.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '32.95'
$ws.Range('E26').Value = '  +11.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.14'
$ws.Range('E27').Value = '  +13.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.79'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.61'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.75'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.03'
$ws.Range('E31').Value = '  +5.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.170'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.74'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0500'
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.94'
$ws.Range('E37').Value = '  +4.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.37'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.134'
$ws.Range('E40').Value = '  +7.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.317'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.32'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.25'
$ws.Range('E44').Value = '  +8.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.55'
$ws.Range('E45').Value = '  +15.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.99'
$ws.Range('E46').Value = '  +1.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.57'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.28'
$ws.Range('E48').Value = '  +5.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.143'
$ws.Range('E49').Value = '  +19.47%  '
$ws.Range('E50').Value = '  +7.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.768.63'
$ws.Range('E51').Value = '  -0.15%  '
